# LOQ4085.xlsx edit script
# Applies the content corrections described in the commit diff:
#  - Fixes "Objetivos:" (B10/C10) to show the actual Portuguese objectives text
#    instead of the teacher name that had been misplaced there.
#  - Inserts a new row (13) to hold the teacher name
#    ("8151869 - Livia Chaguri e Carvalho") under "Docentes responsaveis:".
#  - Fixes "Programa resumido:" / "Programa:" content that had been misaligned.
#  - Fixes "Metodo:" / "Criterio:" / "Norma de recuperacao:" content that had
#    been misaligned (shifted by one row).
#  - Adds the Bibliografia content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos: ---------------------------------------------------
$objetivosPt = "Aplicar os fundamentos teóricos das operações unitárias envolvendo sistemas fluidos e particulados, baseados nos princípios dos fenômenos de transporte I."
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# --- Insert new row 13 for the teacher name (Docentes responsaveis:) ------
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()

$docentes = "8151869 - Livia Chaguri e Carvalho"

# Copy formats from existing column B/C cells so the new row gets the
# correct styles (wrap text, vertical top, red font for column C).
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13").Value = $docentes
$ws.Range("C13").Value = $docentes

# --- Row 14: Programa resumido: (was row 13, now shifted to 14) -----------
$programaResumido = "1)Transporte de fluidos (Newtonianos e não Newtonianos)2)Agitação e mistura3)Caracterização e dinâmica de partículas4)Separação de partículas por ação gravitacional e centrífuga5)Interação sólido – fluido6)Filtração7)Sedimentação"
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# --- Row 16: Programa: (was row 15, now shifted to 16) --------------------
$programaFull = "1)Transporte de fluidos: Tipos de bombas e compressores. Medidores de vazão. Curvas características. Cavitação e altura de sucção disponível (NPSH). Dimensionamento do sistema de bombeamento.2)Agitação e mistura: Tipos de equipamentos e impelidores. Mistura de líquidos. Cálculos de potência de agitadores.3)Caracterização e dinâmica de partículas: Características físicas de partícula isolada. Tamanho de partículas. Peneiramento. Análise granulométrica. Velocidade terminal.4)Separação de partículas por ação gravitacional e centrífuga: Elutriação. Câmara de poeira. Ciclones e centrífugas.5)Interação sólido – fluido: Escoamento em meio poroso. Fluidização.6)Filtração: Tipos de equipamentos. Filtração a pressão e vazão constante. Tortas compressíveis e incompressíveis.7)Sedimentação: Tipos de equipamentos. Cálculo da área e altura de sedimentadores."
$ws.Range("B16").Value = $programaFull
$ws.Range("C16").Value = $programaFull

# --- Row 19: Método: (was row 18, now shifted to 19) -----------------------
$metodo = "Aplicação de 2 provas (P1 e P2)."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Row 20: Critério: (was row 19, now shifted to 20) ---------------------
$criterio = @"
A média do período (MP) será calculada por: MP = (P1+P2)/2. 
Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). 
Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). 
Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental).
"@
$criterio = $criterio.TrimEnd()
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Row 21: Norma de recuperação: (was row 20, now shifted to 21) --------
$normaRecuperacao = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação"
$ws.Range("B21").Value = $normaRecuperacao
$ws.Range("C21").Value = $normaRecuperacao

# --- Row 22: Bibliografia: (was row 21, now shifted to 22) ----------------
$bibliografia = "1)COULSON, J. M.; RICHARDSON; J.F. Chemical Engineering. v.2: Particle Technology e Separation Processes. 5ed. Amsterdan: Butterworth Heinemann, 1229p. 2005;2)COULSON & Richardson's Chemical Engineering: chemical engineering design by R.K. Sinnott. 6ed. Amsterdam: Elsevier Butterworth Heinemann, 895p. 2004;3)COUPER, J. R.; PENNEY, W. R.; FAIR, J. R.; W.; Stanley. M. Chemical Process Equipment: Selection and Design. 2ed. Amsterdam: Elsevier, 814p. 2005;4)MORAES JUNIOR, D. Transporte de líquidos e gases. v.1. São Carlos: Ufscar, 1988;5)FOUST, A. S.; WENZEL, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSEN, L. B. 2ed. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 670p. 2008;6)GEANKOPLIS, C. J. Transport Processes and Separation Process Principles. 4ed. New York: Prentice Hall, 1026p. 2010;7)MCCABE, W. L.; SMITH, J. C.; HARRIOT, P. Unit operations of chemical engineering. 7ed. Boston: McGraw-Hill, 1140 p. 2005;8)PERRY's chemical engineers handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry New York: McGraw-Hill, 2008."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
